$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.399968028068542
$ws.Range("B1").Value = 2.075355291366577
$ws.Range("C1").Value = 3.923242092132568
$ws.Range("D1").Value = 1.138980507850647
$ws.Range("E1").Value = 0.6781141757965088
